$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vehicles")

$ws.Range("A2").Value = "SPRINTER: BT-48331 ; RAM-Promaster 2500 (2021)"
$ws.Range("A3").Value = "SPRINTER: BV-14827 ; MERCEDES- 2500 Cargo Van (2013)"
$ws.Range("A4").Value = "SPRINTER: CA-30264 ; MERCEDES-2500 Cargo Van (2024)"
$ws.Range("A5").Value = "TRUCK: CD-53631 ; ISUZU- NRR (2019)"
